$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link for rows that were reordered (FraxShare, Hedera, InternetComputer(DFINITY))
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

# Update Price (D) and Volume(1h) (E) values for all data rows.
# Price values that look like plain numbers must be forced to remain text
# (matching the original inlineStr cell type) by temporarily applying a text
# number format, then resetting the cell style back to Normal afterwards so
# no stray style/quote-prefix indicator is left behind.
$ws.Range("D2").Value = "27.522.93"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "1.789.68"
$ws.Range("E3").Value = "  +4.32%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.43"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5342"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +11.41%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3767"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +8.32%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "43.03"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.73%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07468"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.29%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.105"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +6.07%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.84"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.01%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.146"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.36%  "
$ws.Range("D15").Value = "1.779.13"
$ws.Range("E15").Value = "  +3.95%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.046"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.06%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "90.24"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.44%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001059"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.17%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06451"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "16.90"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.37%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.936"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.88%  "
$ws.Range("D23").Value = "27.542.72"
$ws.Range("E23").Value = "  +1.93%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.55%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.097"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "156.05"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.47"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.71%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.398"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +15.54%  "
$ws.Range("D29").Value = "1.990.21"
$ws.Range("E29").Value = "  +4.31%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "121.78"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.101"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.28%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.1022"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +11.76%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.624"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +5.54%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.626"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.79%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.02268"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.13%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "8.536"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +14.85%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06014"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.74%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.76%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2074"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.65%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.946"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.81%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6175"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.92%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.420"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.32%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9978"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.143"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.96%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.38"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.41%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5822"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.49%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.630"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "121.59"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.29%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.906"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.12%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.128"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.75%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06744"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
